$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Power draw (W)" column (G) ---

# Header cell G1: same text style as the other header cells (F1)
$ws.Range("G1").Value = "Power draw (W)"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# G2 = B2*D2 (standalone formula, like E2/F2 in row 2)
$ws.Range("G2").Formula = "=B2*D2"

# G3:G40 share a single formula definition, mirroring the existing D/E/F pattern
$ws.Range("G3:G40").Formula = "=B3*D3"

# Writing formulas that reference the (custom-formatted) D column causes the
# engine to copy D's number format onto G; reset G back to the workbook's
# default/general style to match the original column formatting.
$ws.Range("G2:G40").Style = "Normal"

# --- View state updates recorded for this sheet ---
$win = $excel.ActiveWindow
$win.Zoom = 59
$ws.Range("D1").Select()
$ws.Range("AC16").Select()

# --- Move/resize the chart to make room for the new column ---
$co = $ws.ChartObjects().Item(1)
$co.Left = 1203.9485805241143
$co.Top = 83.07606299212598
$co.Width = 1028.0944291338583
$co.Height = 503.26692913385824
